# Option to calculate simple rate instead of compound.
# Adds a new "simple rateloa" worksheet with a simple-interest loan calculator,
# and nudges the selection on the "lenders spread" sheet.

$wb = $excel.ActiveWorkbook

# --- Update selection on the "lenders spread" sheet (B7 -> B6) ---
$wsSpread = $wb.Worksheets.Item("lenders spread")
$wsSpread.Activate() | Out-Null
$wsSpread.Range("B6").Select() | Out-Null

# --- Add the new worksheet as the last tab ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "simple rateloa"

# Column headers
$ws.Range("A1").Value = "loan amount"
$ws.Range("B1").Value = "interest"
$ws.Range("C1").Value = "years"
$ws.Range("D1").Value = "simple interest"
$ws.Range("E1").Value = "total repayment"
$ws.Range("F1").Value = "monthly repayment"

# Row 2 data
$ws.Range("A2").Value = 18000
$ws.Range("B2").Value = 0.06
$ws.Range("C2").Value = 3
$ws.Range("D2").Formula = "=`$A2*`$B2*`$C2"
$ws.Range("E2").Formula = "=`$A2+`$D2"
$ws.Range("F2").Formula = "=`$E2/(`$C2*12)"

# Row 3 data
$ws.Range("A3").Value = 1100
$ws.Range("B3").Value = 0.069
$ws.Range("C3").Value = 3
$ws.Range("D3").Formula = "=`$A3*`$B3*`$C3"
$ws.Range("E3").Formula = "=`$A3+`$D3"
$ws.Range("F3").Formula = "=`$E3/(`$C3*12)"

# Row 4 data
$ws.Range("A4").Value = 1000
$ws.Range("B4").Value = 0.07
$ws.Range("C4").Value = 3
$ws.Range("D4").Formula = "=`$A4*`$B4*`$C4"
$ws.Range("E4").Formula = "=`$A4+`$D4"
$ws.Range("F4").Formula = "=`$E4/(`$C4*12)"

# Column widths (character units, closest achievable to the authored widths)
$ws.Columns.Item(1).ColumnWidth = 10.666666666666666
$ws.Columns.Item(4).ColumnWidth = 14.166666666666666
$ws.Columns.Item(5).ColumnWidth = 13.333333333333334
$ws.Columns.Item(6).ColumnWidth = 14.666666666666666

# Make this the active sheet / tab with the expected selection
$ws.Activate() | Out-Null
$ws.Range("F6").Select() | Out-Null
